# Add a new database entry "test_game_logDB" to the "Database" sheet,
# mirroring the existing "test_game_system" row's layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# Copy formatting from row 2 into row 3 so the new row matches the style
# of the existing data row, then fill in the new values.
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A3:F3").PasteSpecial(-4122) | Out-Null

$ws.Range("A3").Value = "test_game_logDB"
$ws.Range("B3").Value = "localhost"
$ws.Range("C3").Value = "gameAdmin2"
$ws.Range("D3").Value = "admin00!!"
$ws.Range("E3").Value = "game_log"
$ws.Range("F3").Value = "游戏log"

$ws.Rows.Item(3).RowHeight = 13.5

$ws.Range("F4").Select() | Out-Null
